$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-02 Saturday" "2024-03-03 Sunday"

Replace-Text "113×8=904" "562×2=1124"
Replace-Text "429×9=3861" "829×2=1658"
Replace-Text "290×8=2320" "293×5=1465"
Replace-Text "263×4=1052" "752×7=5264"
Replace-Text "111×6=666" "800×2=1600"

Replace-Text "603×7=4221" "573×4=2292"
Replace-Text "707×3=2121" "734×9=6606"
Replace-Text "107×2=214" "541×2=1082"
Replace-Text "786×4=3144" "613×9=5517"
Replace-Text "996×2=1992" "567×6=3402"

Replace-Text "218×5=1090" "576×3=1728"
Replace-Text "212×9=1908" "707×3=2121"
Replace-Text "323×4=1292" "906×7=6342"
Replace-Text "830×5=4150" "716×3=2148"
Replace-Text "920×3=2760" "949×8=7592"

Replace-Text "831×5=4155" "410×7=2870"
Replace-Text "848×9=7632" "820×4=3280"
Replace-Text "467×2=934" "740×4=2960"
Replace-Text "892×4=3568" "418×5=2090"
Replace-Text "946×9=8514" "631×9=5679"

Replace-Text "358×6=2148" "936×7=6552"
Replace-Text "644×4=2576" "401×7=2807"
Replace-Text "133×2=266" "891×5=4455"
Replace-Text "584×7=4088" "871×2=1742"
Replace-Text "741×2=1482" "119×6=714"

Write-Host "All replacements applied"
